# Apply the latest cryptos-list scrape: refreshed prices / 1h volume %
# for most rows, plus two coins (USDe/EthereumClassic and Mantle/Filecoin)
# that swapped rank position between rows 33-34 and 45-46 respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.904.85"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "'3.383.07"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'562.66"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "'154.57"
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'3.384.29"
$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("D9").Value = "'0.542"
$ws.Range("E9").Value = "  +2.15%  "

$ws.Range("D10").Value = "'7.37"
$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").Value = "'3.968.15"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("E14").Value = "  -3.66%  "

$ws.Range("E15").Value = "  +3.05%  "

$ws.Range("D16").Value = "'26.87"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").Value = "'63.018.72"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'3.353.75"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("E19").Value = "  -4.50%  "

$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").Value = "'376.17"
$ws.Range("E21").Value = "  -3.52%  "

$ws.Range("E22").Value = "  -4.86%  "

$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").Value = "'71.31"
$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("E25").Value = "  -2.87%  "

$ws.Range("E26").Value = "  +19.93%  "

$ws.Range("D27").Value = "'9.42"
$ws.Range("E27").Value = "  +6.24%  "

$ws.Range("E28").Value = "  -2.77%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").Value = "'5.99"
$ws.Range("E30").Value = "  +5.48%  "

$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").Value = "'1.32"
$ws.Range("E32").Value = "  +0.87%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'22.95"
$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'6.34"
$ws.Range("E35").Value = "  -5.09%  "

$ws.Range("D36").Value = "'6.72"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").Value = "'157.67"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("E38").Value = "  -3.14%  "

$ws.Range("D39").Value = "'0.0757"
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("D40").Value = "'2.894.69"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("E41").Value = "  -4.32%  "

$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("D43").Value = "'0.0315"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.31"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.751"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("D47").Value = "'23.06"
$ws.Range("E47").Value = "  +3.73%  "

$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "'2.13"
$ws.Range("E49").Value = "  +16.22%  "

$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("E51").Value = "  +2.57%  "
